$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the match data (columns F..V) between row 22 and row 23 ---
# (columns A..E — Indice/pais/torneio/temporada/data_partida — stay put)
$row22vals = @()
$row23vals = @()
for ($col = 6; $col -le 22; $col++) {
    $row22vals += ,$ws.Cells.Item(22, $col).Value2
    $row23vals += ,$ws.Cells.Item(23, $col).Value2
}
for ($i = 0; $i -lt 17; $i++) {
    $col = 6 + $i
    $ws.Cells.Item(22, $col).Value = $row23vals[$i]
    $ws.Cells.Item(23, $col).Value = $row22vals[$i]
}

# --- 2) Append 3 new match rows (103, 104, 105) ---
$newRows = @(
    @{ Row=103; Indice=102; Data=45235.66666666666; F="St. Liege"; G=1; H="KV Mechelen"; I=1;
       J=1.91; K="29/10/2023 18:42"; L=1.8; M="05/11/2023 15:52";
       N=3.83; O="29/10/2023 18:42"; P=3.88; Q="05/11/2023 15:52";
       R=3.58; S="29/10/2023 18:42"; T=4.51; U="05/11/2023 15:59";
       V="https://www.betexplorer.com/football/belgium/jupiler-pro-league/st-liege-kv-mechelen/rRSzKrCc/" },
    @{ Row=104; Indice=103; Data=45235.77083333334; F="Royale Union SG"; G=2; H="Club Brugge KV"; I=1;
       J=2.23; K="29/10/2023 18:42"; L=2.26; M="05/11/2023 18:22";
       N=3.7; O="29/10/2023 18:42"; P=3.61; Q="05/11/2023 18:27";
       R=3.07; S="29/10/2023 18:42"; T=3.18; U="05/11/2023 18:22";
       V="https://www.betexplorer.com/football/belgium/jupiler-pro-league/royale-union-sg-club-brugge/OxpYI2R3/" },
    @{ Row=105; Indice=104; Data=45235.80208333334; F="Charleroi"; G=1; H="Gent"; I=3;
       J=3.26; K="29/10/2023 19:43"; L=3.56; M="05/11/2023 19:13";
       N=3.69; O="29/10/2023 19:43"; P=3.81; Q="05/11/2023 19:13";
       R=2.07; S="29/10/2023 19:43"; T=2.03; U="05/11/2023 19:13";
       V="https://www.betexplorer.com/football/belgium/jupiler-pro-league/charleroi-gent/djrtHtdG/" }
)

foreach ($rec in $newRows) {
    $r = $rec.Row

    # Clone formatting (including style indices) from the row above, so the
    # new rows look identical to the rest of the table (bold/bordered index
    # column, datetime-formatted data_partida column).
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
    $ws.Range("E" + ($r - 1)).Copy()
    $ws.Range("E" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $rec.Indice
    $ws.Cells.Item($r, 2).Value = "belgium"
    $ws.Cells.Item($r, 3).Value = "jupiler-pro-league"
    $ws.Cells.Item($r, 4).Value = "2023-2024"
    $ws.Cells.Item($r, 5).Value = $rec.Data
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
    $ws.Cells.Item($r, 8).Value = $rec.H
    $ws.Cells.Item($r, 9).Value = $rec.I
    $ws.Cells.Item($r, 10).Value = $rec.J
    $ws.Cells.Item($r, 11).Value = $rec.K
    $ws.Cells.Item($r, 12).Value = $rec.L
    $ws.Cells.Item($r, 13).Value = $rec.M
    $ws.Cells.Item($r, 14).Value = $rec.N
    $ws.Cells.Item($r, 15).Value = $rec.O
    $ws.Cells.Item($r, 16).Value = $rec.P
    $ws.Cells.Item($r, 17).Value = $rec.Q
    $ws.Cells.Item($r, 18).Value = $rec.R
    $ws.Cells.Item($r, 19).Value = $rec.S
    $ws.Cells.Item($r, 20).Value = $rec.T
    $ws.Cells.Item($r, 21).Value = $rec.U
    $ws.Cells.Item($r, 22).Value = $rec.V
}
